$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Distribution")

# Update descriptions first (R2: 330Ω->220Ω, R3: 10Ω->22Ω)
$ws.Range("C7").Value = "220Ω, class F"
$ws.Range("C8").Value = "22Ω, class F"

# Then update part numbers to match
$ws.Range("H7").Value = "10-ERJ-U03F2200VCT-ND"
$ws.Range("H8").Value = "10-ERJ-U03F22R0VCT-ND"

# Clear the stale cell selection stored in the sheet view
$ws.Range("A1").Select()
